$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds ISO-formatted dates as plain text (e.g. "2025-09-16").
# Assigning that string directly would make Excel auto-convert it into a
# date serial number, so force the cell to text first.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2025-09-16"

$ws.Range("B22").Value = "21:22:01"
$ws.Range("C22").Value = "1.00 EUR = 1,745.0459"
